$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1702.25
$ws.Range("J2").Value = 2074.6667
$ws.Range("L2").Value = 2074.6667
$ws.Range("N2").Value = -2300.6667

$ws.Range("H28").Value = 7450.8184
$ws.Range("J28").Value = 14698.6
$ws.Range("L28").Value = 14698.6
$ws.Range("N28").Value = -15668.6

$ws.Range("H43").Value = 4003.25
$ws.Range("J43").Value = 4043.5
$ws.Range("L43").Value = 4043.5
$ws.Range("N43").Value = -4181.5

$ws.Range("H64").Value = 4501
$ws.Range("I64").Value = 4501
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 4501
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -4253
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 4501
$ws.Range("I67").Value = 4501
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 4501
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -3643
$ws.Range("N67").ClearContents()

$ws.Range("H70").Value = 1527501.9
$ws.Range("J70").Value = 3649
$ws.Range("L70").Value = 10947
$ws.Range("N70").Value = -11487

$ws.Range("H73").Value = 1527501.9
$ws.Range("J73").Value = 3649
$ws.Range("L73").Value = 10947
$ws.Range("N73").Value = -12819

$ws.Range("H104").Value = 147.4
$ws.Range("J104").Value = 225
$ws.Range("L104").Value = 675
$ws.Range("N104").Value = -4169

$ws.Range("H125").Value = 5817.222
$ws.Range("I125").Value = 812
$ws.Range("K125").Value = 7308
$ws.Range("M125").Value = -4848

$ws.Range("H138").Value = 2623.48
$ws.Range("J138").Value = 4243.6943
$ws.Range("L138").Value = 12731.0829
$ws.Range("N138").Value = -23011.0829

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6952.857
$ws.Range("I45").Value = 4990
$ws.Range("J45").Value = 7280
$ws.Range("K45").Value = 4990
$ws.Range("L45").Value = 7280
$ws.Range("M45").Value = -4613
$ws.Range("N45").Value = -8034

$ws.Range("H63").Value = 3462.2856
$ws.Range("I63").Value = 3420.923
$ws.Range("J63").Value = 4000
$ws.Range("K63").Value = 3420.923
$ws.Range("L63").Value = 4000
$ws.Range("M63").Value = -2734.923
$ws.Range("N63").Value = -5372

$ws.Range("H66").Value = 3462.2856
$ws.Range("I66").Value = 3420.923
$ws.Range("J66").Value = 4000
$ws.Range("K66").Value = 17104.615
$ws.Range("L66").Value = 20000
$ws.Range("M66").Value = -13672.615
$ws.Range("N66").Value = -26864

$ws.Range("H74").Value = 557566
$ws.Range("I74").Value = 626694.3
$ws.Range("K74").Value = 626694.3
$ws.Range("M74").Value = -625820.3

$ws.Range("H77").Value = 557566
$ws.Range("I77").Value = 626694.3
$ws.Range("K77").Value = 3133471.5
$ws.Range("M77").Value = -3129103.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 11113619
$ws.Range("I134").Value = 2315.8
$ws.Range("K134").Value = 6947.400000000001
$ws.Range("M134").Value = -4412.400000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H47").Value = 1000064
$ws.Range("I47").Value = 1000064
$ws.Range("K47").Value = 1000064
$ws.Range("M47").Value = -999498

$ws.Range("H58").Value = 3023.3572
$ws.Range("I58").Value = 3173.4285
$ws.Range("K58").Value = 3173.4285
$ws.Range("M58").Value = -2970.4285

$ws.Range("H62").Value = 14525.6
$ws.Range("I62").Value = 16989
$ws.Range("K62").Value = 16989
$ws.Range("M62").Value = -16365

$ws.Range("H65").Value = 14525.6
$ws.Range("I65").Value = 16989
$ws.Range("K65").Value = 84945
$ws.Range("M65").Value = -81825

$ws.Range("H132").Value = 4115.7
$ws.Range("I132").Value = 2726.3333
$ws.Range("J132").Value = 6199.75
$ws.Range("K132").Value = 8178.999899999999
$ws.Range("L132").Value = 18599.25
$ws.Range("M132").Value = -5648.999899999999
$ws.Range("N132").Value = -23659.25

$ws.Range("H134").Value = 3770.4546
$ws.Range("I134").Value = 3727.5
$ws.Range("K134").Value = 11182.5
$ws.Range("M134").Value = -8647.5

$ws.Range("H136").Value = 3023.3572
$ws.Range("I136").Value = 3173.4285
$ws.Range("K136").Value = 9520.2855
$ws.Range("M136").Value = -6970.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 181.4
$ws.Range("I18").Value = 181.4
$ws.Range("K18").Value = 544.2
$ws.Range("M18").Value = -375.2

$ws.Range("H130").Value = 13283.667
$ws.Range("J130").Value = 14887.167
$ws.Range("L130").Value = 44661.501
$ws.Range("N130").Value = -54701.501

$ws.Range("H131").Value = 4674.227
$ws.Range("J131").Value = 5131.353
$ws.Range("L131").Value = 15394.059
$ws.Range("N131").Value = -25474.059

$ws.Range("H140").Value = 9195.125
$ws.Range("I140").Value = 5382.25
$ws.Range("K140").Value = 16146.75
$ws.Range("M140").Value = -10966.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws.Range("H80").Value = 3287.6155
$ws.Range("I80").Value = 2166.25
$ws.Range("J80").Value = 5081.8
$ws.Range("K80").Value = 2166.25
$ws.Range("L80").Value = 5081.8
$ws.Range("M80").Value = -1168.25
$ws.Range("N80").Value = -7077.8

$ws.Range("H83").Value = 3287.6155
$ws.Range("I83").Value = 2166.25
$ws.Range("J83").Value = 5081.8
$ws.Range("K83").Value = 10831.25
$ws.Range("L83").Value = 25409
$ws.Range("M83").Value = -5839.25
$ws.Range("N83").Value = -35393

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7977.933
$ws.Range("I7").Value = 8097.1816
$ws.Range("J7").Value = 7650
$ws.Range("K7").Value = 8097.1816
$ws.Range("L7").Value = 7650
$ws.Range("M7").Value = -7985.1816
$ws.Range("N7").Value = -7874

$ws.Range("H82").Value = 4297.3184
$ws.Range("J82").Value = 4591.077
$ws.Range("L82").Value = 4591.077
$ws.Range("N82").Value = -5313.077

$ws.Range("H85").Value = 4297.3184
$ws.Range("J85").Value = 4591.077
$ws.Range("L85").Value = 4591.077
$ws.Range("N85").Value = -7087.077

$ws.Range("H126").Value = 7977.933
$ws.Range("I126").Value = 8097.1816
$ws.Range("J126").Value = 7650
$ws.Range("K126").Value = 24291.5448
$ws.Range("L126").Value = 22950
$ws.Range("M126").Value = -21821.5448
$ws.Range("N126").Value = -27890

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 14428.571
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()

$ws.Range("H37").Value = 89497.5
$ws.Range("J37").Value = 79995
$ws.Range("L37").Value = 79995
$ws.Range("N37").Value = -80401

$ws.Range("H81").Value = 1511.7693
$ws.Range("J81").Value = 1998
$ws.Range("L81").Value = 3996
$ws.Range("N81").Value = -6118

$ws.Range("H84").Value = 1511.7693
$ws.Range("J84").Value = 1998
$ws.Range("L84").Value = 19980
$ws.Range("N84").Value = -30588

$ws.Range("H132").Value = 205160.4
$ws.Range("I132").Value = 1065.762
$ws.Range("J132").Value = 1429728.2
$ws.Range("K132").Value = 3197.286
$ws.Range("L132").Value = 4289184.6
$ws.Range("M132").Value = -667.2860000000001
$ws.Range("N132").Value = -4294244.6

